# Generate Report for Handback
# Update the "last generated/handoff/handback" timestamps recorded on the
# Overview, zh-cn and de-de sheets to reflect a fresh report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview.Range("G2").Value = "2016-08-27 07:02:17"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the first file.
$wsZhCn.Range("H2").Value = "2016-08-27 07:02:12"
$wsZhCn.Range("K2").Value = "2016-08-27 07:02:29"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the first file.
$wsDeDe.Range("H2").Value = "2016-08-27 07:02:17"
$wsDeDe.Range("K2").Value = "2016-08-27 07:02:36"
